$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, [string]$Text)
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "303.02"
Set-TextValue $ws.Range("E2") "5.32%"
Set-TextValue $ws.Range("D3") "34.76"
Set-TextValue $ws.Range("E3") "12.26%"
Set-TextValue $ws.Range("D4") "5.170"
Set-TextValue $ws.Range("E4") "4.80%"
Set-TextValue $ws.Range("D5") "0.07810"
Set-TextValue $ws.Range("E5") "6.70%"
Set-TextValue $ws.Range("D6") "2.332"
Set-TextValue $ws.Range("E6") "-2.79%"
Set-TextValue $ws.Range("D7") "8.009"
Set-TextValue $ws.Range("E7") "3.65%"
Set-TextValue $ws.Range("D8") "3.990"
Set-TextValue $ws.Range("E8") "7.15%"
Set-TextValue $ws.Range("D9") "0.9302"
Set-TextValue $ws.Range("E9") "2.95%"
Set-TextValue $ws.Range("D10") "0.1017"
Set-TextValue $ws.Range("E10") "8.93%"
Set-TextValue $ws.Range("D11") "0.1836"
Set-TextValue $ws.Range("E11") "8.90%"
Set-TextValue $ws.Range("D12") "0.08497"
Set-TextValue $ws.Range("E12") "4.07%"
Set-TextValue $ws.Range("D13") "0.03486"
Set-TextValue $ws.Range("E13") "11.56%"
Set-TextValue $ws.Range("D14") "0.09908"
Set-TextValue $ws.Range("E14") "-0.21%"
Set-TextValue $ws.Range("D15") "0.001477"
Set-TextValue $ws.Range("E15") "-1.44%"
Set-TextValue $ws.Range("D16") "0.04602"
Set-TextValue $ws.Range("E16") "2.25%"
Set-TextValue $ws.Range("D17") "0.005780"
Set-TextValue $ws.Range("E17") "0.63%"
Set-TextValue $ws.Range("D18") "3.474"
Set-TextValue $ws.Range("E18") "-0.61%"
Set-TextValue $ws.Range("E19") "0.36%"
Set-TextValue $ws.Range("E20") "3.41%"
Set-TextValue $ws.Range("E21") "-0.36%"
Set-TextValue $ws.Range("D22") "4.543"
Set-TextValue $ws.Range("E22") "7.78%"
Set-TextValue $ws.Range("E24") "0.79%"
Set-TextValue $ws.Range("D25") "0.004439"
Set-TextValue $ws.Range("E25") "6.74%"
Set-TextValue $ws.Range("E26") "0.04%"
Set-TextValue $ws.Range("D27") "0.0003398"
Set-TextValue $ws.Range("E27") "0.12%"
Set-TextValue $ws.Range("D39") "0.01759"
Set-TextValue $ws.Range("E39") "11.64%"
Set-TextValue $ws.Range("D40") "0.04721"
Set-TextValue $ws.Range("E40") "6.16%"
Set-TextValue $ws.Range("D41") "0.007678"
Set-TextValue $ws.Range("E41") "4.09%"
Set-TextValue $ws.Range("D42") "0.1406"
Set-TextValue $ws.Range("E42") "6.13%"
Set-TextValue $ws.Range("D43") "0.007039"
Set-TextValue $ws.Range("E43") "-25.33%"
Set-TextValue $ws.Range("D44") "0.002280"
Set-TextValue $ws.Range("E44") "1.83%"
Set-TextValue $ws.Range("D45") "0.009334"
Set-TextValue $ws.Range("E45") "4.19%"
Set-TextValue $ws.Range("D46") "0.00005992"
Set-TextValue $ws.Range("E46") "-1.85%"
Set-TextValue $ws.Range("E48") "8.97%"
Set-TextValue $ws.Range("D49") "0.002696"
Set-TextValue $ws.Range("E49") "34.72%"
